$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Data")
$wsQuestions = $wb.Worksheets.Item("Questions")

# --- Data sheet: move the selection/active cell (view-only change) ---
$wsData.Range("C23").Select()

# --- Questions sheet ---
# Clear the old long "Retour blended" text from A12 but keep its cell style
$wsQuestions.Range("A12").ClearContents()

# The text got split into two separate cells further down the sheet
$wsQuestions.Range("A14").Value = "BDMOYENNE"
$wsQuestions.Range("A15").Value = "DAVERAGE"

# Re-activate Questions (it was the active tab) and move its selection too
$wsQuestions.Activate()
$wsQuestions.Range("A12").Select()
